# Tripadvisor New Orleans shard_63 update:
#   1. Reorder worksheets so "review_info" comes before "hotel_info".
#   2. Add a new "State" column to "hotel_info" right after "Hotel_Name",
#      populated with "Louisiana" for the existing hotel row.

$wb = $excel.ActiveWorkbook

$hotel  = $wb.Worksheets.Item("hotel_info")
$review = $wb.Worksheets.Item("review_info")

# Move review_info so it sits immediately before hotel_info in the tab order.
$review.Move($hotel)

# Sheet objects here are resolved by slot/position, so after the Move the
# variable above that used to point at "hotel_info" now reports the other
# sheet's name. Re-fetch the worksheet by name to get a fresh, correct
# reference to work with.
$hotel = $wb.Worksheets.Item("hotel_info")

# Insert a new blank column C (shifting City/Zip/etc. one column right) and
# fill in the State header + value.
$hotel.Columns.Item(3).Insert()
$hotel.Cells.Item(1, 3).Value = "State"
$hotel.Cells.Item(2, 3).Value = "Louisiana"
